$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the extended-cost column (I) for every BOM line: qty * unit price.
# Row 2 carries the "master" formula of the shared-formula group I2:I37,
# the rest just reference it (mirrors what Excel does on a fill-down).
$ws.Range("I2:I37").Formula = '=$E2*G2'

# Grand-total row underneath the table.
$ws.Range("I38").Formula = '=SUM(I2:I37)'

# Update the on-screen selection/view to match the new working area,
# and drop the old scrolled-in topLeftCell.
[void]$ws.Range("E1:H38").Select()
